$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Brandon"
$ws.Range("B15").Value = "Provost"
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 12
$ws.Range("E15").Value = 2003

$ws.Range("A16").Value = "Ian"
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = 2003

$ws.Range("A17").Value = "Edmond"
$ws.Range("B17").Value = "Simonian"
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = 2004

$ws.Range("E22").Select()
